$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the old "DOI" (C) and "SpeciesAuthor" (D) columns - the
#    remaining columns shift two places to the left.
# ------------------------------------------------------------------
$ws.Columns("C:D").Delete()

# ------------------------------------------------------------------
# 2. Insert four new blank columns right after "AdditionalSources"
#    (now column G) to hold the corrected-citation fields.
# ------------------------------------------------------------------
$ws.Columns("H:K").Insert()

# ------------------------------------------------------------------
# 3. Populate the headers for the new columns and give them the
#    same bold header style the rest of row 1 uses.
# ------------------------------------------------------------------
$ws.Range("H1").Value = "Author.corr"
$ws.Range("I1").Value = "Journal.corr"
$ws.Range("J1").Value = "YearPublication.corr"
$ws.Range("K1").Value = "DOI.ISBN.corr"
$ws.Range("H1:K1").Font.Bold = $true
$ws.Range("H1:K1").ColumnWidth = 15

# ------------------------------------------------------------------
# 4. Tidy up the view: selection moves to H13, and the frozen
#    top-left cell override is cleared.
# ------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("H13").Select()
